$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 12).Value = 0.9929783193494216
$ws.Cells.Item(21, 5).Value = 0.9882828385668253
$ws.Cells.Item(21, 12).Value = 0.9920501090198105
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 12).Value = 0.9944092447426414

$wb.Save()
